# The sheet currently has:
#   L:P = "Inicio estación más cercana 1..5" (and matching date values)
#   Q:R = "Estación más cercana 6/7"        (and matching station-name values)
#   S:T = "Inicio estación más cercana 6/7" (and matching date values)
#
# The target layout moves "Estación más cercana 6/7" (currently Q:R) to sit
# right after column K (i.e. before the "Inicio..." columns), pushing the
# "Inicio estación más cercana 1..5" block from L:P to N:R. Columns S:T stay
# put logically (still hold "Inicio ... 6/7"), they are just reached via a
# different physical path during the edit.
#
# Implemented as: insert two blank columns at L:M (shifting L: onward right
# by two, so the old Q:R data is now at S:T), copy that data into the new
# L:M, then delete the now-duplicated S:T columns (shifting the remaining
# data, old U:V i.e. "Inicio ... 6/7", back left into S:T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftToRight = -4161
$xlShiftToLeft = -4159

# Step 1: insert two blank columns before column L (shifts L:T -> N:V)
$ws.Range("L1:M40").Insert($xlShiftToRight)

# Step 2: the old "Estación más cercana 6/7" data (previously at Q:R) is now
# at S:T after the shift; copy it into the newly inserted L:M columns.
$ws.Range("L1:M40").Value = $ws.Range("S1:T40").Value()

# Step 3: remove the duplicated S:T columns, shifting the remaining columns
# (old "Inicio estación más cercana 6/7", at U:V) left into S:T.
$ws.Range("S1:T40").Delete($xlShiftToLeft)
